# Auto-generated edit script applying scheduled market-data refresh values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 799.1111
$ws.Range("J32").Value = 828.6
$ws.Range("L32").Value = 828.6
$ws.Range("N32").Value = -1480.6
$ws.Range("H33").Value = 206.1579
$ws.Range("I33").Value = 74.82353000000001
$ws.Range("K33").Value = 74.82353000000001
$ws.Range("M33").Value = 154.17647
$ws.Range("H53").Value = 236.21053
$ws.Range("I53").Value = 229
$ws.Range("J53").Value = 248.57143
$ws.Range("K53").Value = 229
$ws.Range("L53").Value = 248.57143
$ws.Range("M53").Value = 408
$ws.Range("N53").Value = -1522.57143
$ws.Range("H55").Value = 497.73685
$ws.Range("I55").Value = 192.5
$ws.Range("K55").Value = 192.5
$ws.Range("M55").Value = 21.5
$ws.Range("H64").Value = 4942.857
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H67").Value = 4942.857
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H105").Value = 67440
$ws.Range("J105").Value = 67440
$ws.Range("L105").Value = 67440
$ws.Range("N105").Value = -74428
$ws.Range("H129").Value = 2816.1177
$ws.Range("J129").Value = 3462.818
$ws.Range("L129").Value = 10388.454
$ws.Range("N129").Value = -20388.454
$ws.Range("H131").Value = 7912.125
$ws.Range("I131").Value = 1648.5
$ws.Range("K131").Value = 4945.5
$ws.Range("M131").Value = 94.5
$ws.Range("H135").Value = 722.2083
$ws.Range("I135").Value = 506.7
$ws.Range("J135").Value = 1799.75
$ws.Range("K135").Value = 4560.3
$ws.Range("L135").Value = 16197.75
$ws.Range("M135").Value = -2025.3
$ws.Range("N135").Value = -21267.75
$ws.Range("H138").Value = 17777.111
$ws.Range("J138").Value = 19139.6
$ws.Range("L138").Value = 57418.8
$ws.Range("N138").Value = -67698.79999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1538.9166
$ws.Range("I61").Value = 1538.9166
$ws.Range("K61").Value = 1538.9166
$ws.Range("M61").Value = -1326.9166
$ws.Range("H74").Value = 4513.8667
$ws.Range("I74").Value = 1899.4286
$ws.Range("K74").Value = 1899.4286
$ws.Range("M74").Value = -1025.4286
$ws.Range("H77").Value = 4513.8667
$ws.Range("I77").Value = 1899.4286
$ws.Range("K77").Value = 9497.143
$ws.Range("M77").Value = -5129.143
$ws.Range("H122").Value = 347402.9
$ws.Range("I122").Value = 477990.84
$ws.Range("J122").Value = 4609.5
$ws.Range("K122").Value = 1433972.52
$ws.Range("L122").Value = 13828.5
$ws.Range("M122").Value = -1431522.52
$ws.Range("N122").Value = -18728.5
$ws.Range("H132").Value = 2216.5405
$ws.Range("I132").Value = 2167
$ws.Range("K132").Value = 6501
$ws.Range("M132").Value = -3971
$ws.Range("H136").Value = 1538.9166
$ws.Range("I136").Value = 1538.9166
$ws.Range("K136").Value = 4616.7498
$ws.Range("M136").Value = -2066.7498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 873.4
$ws.Range("I134").Value = 603.6316
$ws.Range("K134").Value = 1810.8948
$ws.Range("M134").Value = 724.1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2072.6667
$ws.Range("I16").Value = 2219
$ws.Range("K16").Value = 2219
$ws.Range("M16").Value = -1932
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 200
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -900
$ws.Range("H99").Value = 14279.1
$ws.Range("I99").Value = 12268.2
$ws.Range("K99").Value = 12268.2
$ws.Range("M99").Value = -10770.2
$ws.Range("H107").Value = 905.75
$ws.Range("I107").Value = 473.22223
$ws.Range("K107").Value = 473.22223
$ws.Range("M107").Value = 1446.77777
$ws.Range("H113").Value = 2072.6667
$ws.Range("I113").Value = 2219
$ws.Range("K113").Value = 2219
$ws.Range("M113").Value = -49
$ws.Range("H126").Value = 14279.1
$ws.Range("I126").Value = 12268.2
$ws.Range("K126").Value = 36804.60000000001
$ws.Range("M126").Value = -34334.60000000001
$ws.Range("H134").Value = 4583.5264
$ws.Range("I134").Value = 3313.1538
$ws.Range("K134").Value = 9939.4614
$ws.Range("M134").Value = -7404.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1300.2
$ws.Range("J107").Value = 1437.7
$ws.Range("L107").Value = 4313.1
$ws.Range("N107").Value = -8153.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 7048
$ws.Range("I41").Value = 7048
$ws.Range("K41").Value = 7048
$ws.Range("M41").Value = -6693
$ws.Range("H122").Value = 1002870.25
$ws.Range("I122").Value = 145083.42
$ws.Range("J122").Value = 2503997.2
$ws.Range("K122").Value = 435250.26
$ws.Range("L122").Value = 7511991.600000001
$ws.Range("M122").Value = -432800.26
$ws.Range("N122").Value = -7516891.600000001
$ws.Range("H126").Value = 4828.143
$ws.Range("I126").Value = 4399
$ws.Range("K126").Value = 13197
$ws.Range("M126").Value = -10727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3933.1667
$ws.Range("I7").Value = 3719.8
$ws.Range("K7").Value = 3719.8
$ws.Range("M7").Value = -3607.8
$ws.Range("H22").Value = 932.1429000000001
$ws.Range("I22").Value = 949.125
$ws.Range("J22").Value = 921.6923
$ws.Range("K22").Value = 949.125
$ws.Range("L22").Value = 921.6923
$ws.Range("M22").Value = -654.125
$ws.Range("N22").Value = -1511.6923
$ws.Range("H27").Value = 932.1429000000001
$ws.Range("I27").Value = 949.125
$ws.Range("J27").Value = 921.6923
$ws.Range("K27").Value = 949.125
$ws.Range("L27").Value = 921.6923
$ws.Range("M27").Value = -842.125
$ws.Range("N27").Value = -1135.6923
$ws.Range("H40").Value = 4640.75
$ws.Range("I40").Value = 4504
$ws.Range("J40").Value = 4777.5
$ws.Range("K40").Value = 4504
$ws.Range("L40").Value = 4777.5
$ws.Range("M40").Value = -4368
$ws.Range("N40").Value = -5049.5
$ws.Range("H82").Value = 1804.7894
$ws.Range("I82").Value = 1821.3077
$ws.Range("K82").Value = 1821.3077
$ws.Range("M82").Value = -1460.3077
$ws.Range("H85").Value = 1804.7894
$ws.Range("I85").Value = 1821.3077
$ws.Range("K85").Value = 1821.3077
$ws.Range("M85").Value = -573.3077000000001
$ws.Range("H100").Value = 4828.5713
$ws.Range("I100").Value = 1200
$ws.Range("K100").Value = 1200
$ws.Range("M100").Value = -659
$ws.Range("H126").Value = 3933.1667
$ws.Range("I126").Value = 3719.8
$ws.Range("K126").Value = 11159.4
$ws.Range("M126").Value = -8689.400000000001
$ws.Range("H132").Value = 4465.483
$ws.Range("I132").Value = 2949.0667
$ws.Range("J132").Value = 6090.2144
$ws.Range("K132").Value = 8847.2001
$ws.Range("L132").Value = 18270.6432
$ws.Range("M132").Value = -6317.2001
$ws.Range("N132").Value = -23330.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1343.7
$ws.Range("J107").Value = 4000
$ws.Range("L107").Value = 12000
$ws.Range("N107").Value = -15840
$ws.Range("H122").Value = 1605.9333
$ws.Range("I122").Value = 1542.0714
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4626.2142
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2176.2142
$ws.Range("N122").Value = -12400
$ws.Range("H136").Value = 41394.31
$ws.Range("I136").Value = 2510.889
$ws.Range("J136").Value = 128882
$ws.Range("K136").Value = 7532.667
$ws.Range("L136").Value = 386646
$ws.Range("M136").Value = -4982.667
$ws.Range("N136").Value = -391746
